# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Castle Brite Damasco, Region de O'Higgins)
# at the top of the data block (row 183), pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 183 (existing rows 183-210 shift down to 185-212)
$ws.Rows("183:184").Insert()

# --- New row 183 ---
$ws.Cells.Item(183, 1).Value = 9
$ws.Cells.Item(183, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(183, 3).Value = "Metropolitana"
$ws.Cells.Item(183, 4).Value = 45275
$ws.Cells.Item(183, 5).Value = 13
$ws.Cells.Item(183, 6).Value = "Fruta"
$ws.Cells.Item(183, 7).Value = 100103
$ws.Cells.Item(183, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(183, 9).Value = 100103003
$ws.Cells.Item(183, 10).Value = "Damasco"
$ws.Cells.Item(183, 11).Value = "Castle Brite"
$ws.Cells.Item(183, 12).Value = "Primera"
$ws.Cells.Item(183, 13).Value = 120
$ws.Cells.Item(183, 14).Value = 9000
$ws.Cells.Item(183, 15).Value = 9000
$ws.Cells.Item(183, 16).Value = 9000
$ws.Cells.Item(183, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(183, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(183, 19).Value = 900
$ws.Cells.Item(183, 20).Value = 10

# --- New row 184 ---
$ws.Cells.Item(184, 1).Value = 9
$ws.Cells.Item(184, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(184, 3).Value = "Metropolitana"
$ws.Cells.Item(184, 4).Value = 45275
$ws.Cells.Item(184, 5).Value = 13
$ws.Cells.Item(184, 6).Value = "Fruta"
$ws.Cells.Item(184, 7).Value = 100103
$ws.Cells.Item(184, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(184, 9).Value = 100103003
$ws.Cells.Item(184, 10).Value = "Damasco"
$ws.Cells.Item(184, 11).Value = "Castle Brite"
$ws.Cells.Item(184, 12).Value = "Segunda"
$ws.Cells.Item(184, 13).Value = 180
$ws.Cells.Item(184, 14).Value = 7000
$ws.Cells.Item(184, 15).Value = 7000
$ws.Cells.Item(184, 16).Value = 7000
$ws.Cells.Item(184, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(184, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(184, 19).Value = 700
$ws.Cells.Item(184, 20).Value = 10
